# Update SA license default
# - Change the "use and reproduction" accessCondition on the Template sheet
#   (column BJ, row 1) from the CC BY-NC-SA 3.0 URL to the CC BY-NC-SA 4.0 URL.
# - Change the short license-name text on the Guide sheet (B22) from the old
#   long-form CC 3.0 description to the new short CC 4.0 label.
# - Update the remembered selections on both sheets (Template -> B1, Guide -> B22),
#   leaving the Template sheet as the active sheet/tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Template"
$ws2 = $wb.Worksheets.Item(2)   # "Guide"

# Update the license accessCondition tag used in the Template header row.
$ws1.Range("BJ1").Value = '<mods:accessCondition displayLabel="License" type="use and reproduction" xlink:href="https://creativecommons.org/licenses/by-nc-sa/4.0/">'

# Update the short license text shown in the Guide sheet.
$ws2.Range("B22").Value = '"Attribution-NonCommercial-ShareAlike 4.0 International (CC BY-NC-SA 4.0)"'

# Set the Guide sheet's remembered selection to B22 first...
$ws2.Range("B22").Select() | Out-Null

# ...then reactivate the Template sheet and select B1, so it remains the
# active/visible tab with its own selection remembered.
$ws1.Activate() | Out-Null
$ws1.Range("B1").Select() | Out-Null
